$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cValues = @(380,11,75,1099,146,48,16,116,1452,524,124,145,169,22,15,50,210,2690,2795,991,9,7,223,59,212,97,61,626,155,202,46,51,17,652,643,5,88,387,1422,413,269,344,18,2050,139,209,78,197,418,211,40,87,104,208,8,244,57,69,52,36,14,43,344.3225806451613)
$dValues = @(426.5,9.5,71,1142.5,161.5,56.5,15.5,150,1445,524,138.5,165,218,20.5,13.5,59,210.5,2746.5,2929,1025.5,9,8.5,225,66,214,123,68,632.5,180,208,51,60,15,630,607,26,98.5,363.5,1436,449.5,148,344,16,1811,140.5,218.5,83.5,199.5,406.5,213,37,70.5,107,217,23,244,68,81.5,59,29.5,10.5,53)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}
